$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# The log times in column D were corrected: the new clock time is
# ~0.0255 earlier (11:39:54 -> 11:03:15). D2 and D3 keep full float
# precision while D4:D75 use the slightly truncated literal that Excel
# produced for the remaining rows.
$ws.Range("D2").Value = 0.46059027777777778
$ws.Range("D3").Value = 0.46059027777777778
$ws.Range("D4:D75").Value = 0.460590277777778

# All of D2:D75 end up sharing one uniform time format (h:mm:ss) with an
# explicit black font color, instead of the old split between the
# builtin time format (D2) and the custom "h:mm:ss;@" format (D3:D75).
$ws.Range("D2:D75").NumberFormat = "h:mm:ss"
$ws.Range("D2:D75").Font.Color = 0

# The selection moved from the header row to the edited time column.
$ws.Range("D2:D75").Select()
